$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set D4:D26 Runmode values to "Y" so all notification test cases run.
$ws.Range("D4:D26").Value = "Y"

# Update the selection to reflect the new active cell/range.
$ws.Activate()
$ws.Range("D2:D26").Select()
